$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.853.95"
$ws.Range("E2").Value = "'  -1.48%  "
$ws.Range("D3").Value = "'3.411.89"
$ws.Range("E3").Value = "'  -0.87%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'570.92"
$ws.Range("E5").Value = "'  -0.48%  "
$ws.Range("D6").Value = "'162.87"
$ws.Range("E6").Value = "'  +2.35%  "
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("D8").Value = "'3.411.77"
$ws.Range("E8").Value = "'  -0.89%  "
$ws.Range("D9").Value = "'0.551"
$ws.Range("E9").Value = "'  -5.20%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "'  +1.16%  "
$ws.Range("E11").Value = "'  -2.12%  "
$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = "'  -4.68%  "
$ws.Range("D13").Value = "'4.006.19"
$ws.Range("E13").Value = "'  -0.68%  "
$ws.Range("E14").Value = "'  +1.13%  "
$ws.Range("D15").Value = "'26.96"
$ws.Range("E15").Value = "'  -2.55%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "'  -5.57%  "
$ws.Range("D17").Value = "'63.957.60"
$ws.Range("E17").Value = "'  -1.32%  "
$ws.Range("D18").Value = "'3.426.02"
$ws.Range("E18").Value = "'  -0.43%  "
$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = "'  -3.59%  "
$ws.Range("D20").Value = "'13.60"
$ws.Range("E20").Value = "'  -2.18%  "
$ws.Range("D21").Value = "'376.71"
$ws.Range("E21").Value = "'  -1.27%  "
$ws.Range("D22").Value = "'7.76"
$ws.Range("E22").Value = "'  -2.68%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "'  +0.11%  "
$ws.Range("D24").Value = "'71.07"
$ws.Range("E24").Value = "'  -1.49%  "
$ws.Range("D25").Value = "'0.515"
$ws.Range("E25").Value = "'  -6.02%  "
$ws.Range("D26").Value = "'0.0000115"
$ws.Range("E26").Value = "'  -3.14%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "'  -3.02%  "
$ws.Range("E28").Value = "'  +0.26%  "
$ws.Range("E29").Value = "'  +0.20%  "
$ws.Range("D30").Value = "'6.09"
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("D31").Value = "'1.39"
$ws.Range("E31").Value = "'  -6.01%  "
$ws.Range("E32").Value = "'  -0.23%  "
$ws.Range("D33").Value = "'22.86"
$ws.Range("E33").Value = "'  -1.69%  "
$ws.Range("D34").Value = "'7.08"
$ws.Range("E34").Value = "'  +1.07%  "
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "'  -4.55%  "
$ws.Range("D36").Value = "'159.96"
$ws.Range("E36").Value = "'  -0.81%  "
$ws.Range("D37").Value = "'0.861"
$ws.Range("E37").Value = "'  +10.61%  "
$ws.Range("E38").Value = "'  -4.39%  "
$ws.Range("D39").Value = "'0.0727"
$ws.Range("E39").Value = "'  -2.54%  "
$ws.Range("D40").Value = "'2.785.97"
$ws.Range("E40").Value = "'  -3.94%  "
$ws.Range("D41").Value = "'25.80"
$ws.Range("E41").Value = "'  -1.61%  "
$ws.Range("D42").Value = "'42.78"
$ws.Range("E42").Value = "'  -0.35%  "
$ws.Range("D43").Value = "'6.43"
$ws.Range("E43").Value = "'  -3.71%  "
$ws.Range("D44").Value = "'26.03"
$ws.Range("E44").Value = "'  +0.23%  "
$ws.Range("D45").Value = "'4.40"
$ws.Range("E45").Value = "'  -3.18%  "
$ws.Range("D46").Value = "'0.0306"
$ws.Range("E46").Value = "'  -3.30%  "
$ws.Range("D47").Value = "'2.43"
$ws.Range("E47").Value = "'  +6.82%  "
$ws.Range("D48").Value = "'329.00"
$ws.Range("E48").Value = "'  +4.02%  "
$ws.Range("D49").Value = "'1.04"
$ws.Range("E49").Value = "'  -3.83%  "
$ws.Range("D50").Value = "'6.30"
$ws.Range("E50").Value = "'  -3.22%  "
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "'  -2.99%  "
